$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 232.66667
$ws.Range("I6").Value = 232.66667
$ws.Range("K6").Value = 698.00001
$ws.Range("M6").Value = -586.00001
$ws.Range("H8").Value = 2483.7144
$ws.Range("I8").Value = 2483.7144
$ws.Range("K8").Value = 7451.1432
$ws.Range("M8").Value = -7312.1432
$ws.Range("H62").Value = 1357.4333
$ws.Range("I62").Value = 1228.5
$ws.Range("J62").Value = 1873.1666
$ws.Range("K62").Value = 1228.5
$ws.Range("L62").Value = 1873.1666
$ws.Range("M62").Value = -604.5
$ws.Range("N62").Value = -3121.1666
$ws.Range("H65").Value = 1357.4333
$ws.Range("I65").Value = 1228.5
$ws.Range("J65").Value = 1873.1666
$ws.Range("K65").Value = 6142.5
$ws.Range("L65").Value = 9365.833000000001
$ws.Range("M65").Value = -3022.5
$ws.Range("N65").Value = -15605.833
$ws.Range("H68").Value = 27813.75
$ws.Range("J68").Value = 27813.75
$ws.Range("L68").Value = 27813.75
$ws.Range("N68").Value = -29311.75
$ws.Range("H71").Value = 27813.75
$ws.Range("J71").Value = 27813.75
$ws.Range("L71").Value = 83441.25
$ws.Range("N71").Value = -90929.25
$ws.Range("H98").Value = 44747.89
$ws.Range("I98").Value = 63395.406
$ws.Range("J98").Value = 2125
$ws.Range("K98").Value = 63395.406
$ws.Range("L98").Value = 2125
$ws.Range("M98").Value = -61897.406
$ws.Range("N98").Value = -5121
$ws.Range("H122").Value = 44747.89
$ws.Range("I122").Value = 63395.406
$ws.Range("J122").Value = 2125
$ws.Range("K122").Value = 190186.218
$ws.Range("L122").Value = 6375
$ws.Range("M122").Value = -187736.218
$ws.Range("N122").Value = -11275

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 1621.25
$ws.Range("I3").Value = 1621.25
$ws.Range("K3").Value = 1621.25
$ws.Range("M3").Value = -1506.25
$ws.Range("H6").Value = 5500
$ws.Range("I6").Value = 1000
$ws.Range("J6").Value = 10000
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = -827
$ws.Range("N6").Value = -10346
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H11").Value = 5500
$ws.Range("J11").Value = 6500
$ws.Range("L11").Value = 6500
$ws.Range("N11").Value = -6788
$ws.Range("H32").Value = 5223.64
$ws.Range("I32").Value = 3465.9312
$ws.Range("J32").Value = 16986.77
$ws.Range("K32").Value = 3465.9312
$ws.Range("L32").Value = 16986.77
$ws.Range("M32").Value = -3178.9312
$ws.Range("N32").Value = -17560.77
$ws.Range("H102").Value = 1764.1666
$ws.Range("I102").Value = 1587.619
$ws.Range("K102").Value = 1587.619
$ws.Range("M102").Value = 34.38100000000009

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 258
$ws.Range("I8").Value = 258
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 258
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -118
$ws.Range("N8").ClearContents()
$ws.Range("H11").Value = 1474
$ws.Range("I11").Value = 140.5
$ws.Range("J11").Value = 2540.8
$ws.Range("K11").Value = 140.5
$ws.Range("L11").Value = 2540.8
$ws.Range("M11").Value = -0.5
$ws.Range("N11").Value = -2820.8
$ws.Range("H12").Value = 8000
$ws.Range("J12").Value = 8000
$ws.Range("L12").Value = 8000
$ws.Range("N12").Value = -8336
$ws.Range("H20").Value = 1803.186
$ws.Range("I20").Value = 1742.9354
$ws.Range("K20").Value = 1742.9354
$ws.Range("M20").Value = -1495.9354

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 4954.2856
$ws.Range("I2").Value = 5613.3335
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 5613.3335
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -5500.3335
$ws.Range("N2").Value = -1226
$ws.Range("H4").Value = 959727.3
$ws.Range("I4").Value = 92181.91
$ws.Range("K4").Value = 92181.91
$ws.Range("M4").Value = -92069.91
$ws.Range("H13").Value = 1903.2
$ws.Range("J13").Value = 4000
$ws.Range("L13").Value = 4000
$ws.Range("N13").Value = -4278

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 399.3
$ws.Range("I10").Value = 120.25
$ws.Range("J10").Value = 1515.5
$ws.Range("K10").Value = 360.75
$ws.Range("L10").Value = 4546.5
$ws.Range("M10").Value = -221.75
$ws.Range("N10").Value = -4824.5
$ws.Range("H35").Value = 1237.5
$ws.Range("J35").Value = 1550
$ws.Range("L35").Value = 4650
$ws.Range("N35").Value = -5226
$ws.Range("H59").Value = 1042.5
$ws.Range("I59").Value = 223.33333
$ws.Range("J59").Value = 3500
$ws.Range("K59").Value = 669.99999
$ws.Range("L59").Value = 10500
$ws.Range("M59").Value = -129.99999
$ws.Range("N59").Value = -11580
$ws.Range("H124").Value = 1291.4286
$ws.Range("I124").Value = 509
$ws.Range("J124").Value = 3247.5
$ws.Range("K124").Value = 1527
$ws.Range("L124").Value = 9742.5
$ws.Range("M124").Value = 3383
$ws.Range("N124").Value = -19562.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2433429
$ws.Range("I3").Value = 5668734.5
$ws.Range("J3").Value = 6950
$ws.Range("K3").Value = 5668734.5
$ws.Range("L3").Value = 6950
$ws.Range("M3").Value = -5668618.5
$ws.Range("N3").Value = -7182
$ws.Range("H11").Value = 81679680
$ws.Range("I11").Value = 10862500
$ws.Range("J11").Value = 110006550
$ws.Range("K11").Value = 10862500
$ws.Range("L11").Value = 110006550
$ws.Range("M11").Value = -10862361
$ws.Range("N11").Value = -110006828
$ws.Range("H34").Value = 17250
$ws.Range("J34").Value = 17250
$ws.Range("L34").Value = 17250
$ws.Range("N34").Value = -17786
$ws.Range("H76").Value = 17250
$ws.Range("J76").Value = 17250
$ws.Range("L76").Value = 17250
$ws.Range("N76").Value = -17880
$ws.Range("H79").Value = 17250
$ws.Range("J79").Value = 17250
$ws.Range("L79").Value = 17250
$ws.Range("N79").Value = -19434
$ws.Range("H126").Value = 1900
$ws.Range("I126").Value = 1575
$ws.Range("J126").Value = 2766.6667
$ws.Range("K126").Value = 4725
$ws.Range("L126").Value = 8300.000100000001
$ws.Range("M126").Value = -2255
$ws.Range("N126").Value = -13240.0001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1431
$ws.Range("I7").Value = 1209.2354
$ws.Range("J7").Value = 1969.5714
$ws.Range("K7").Value = 1209.2354
$ws.Range("L7").Value = 1969.5714
$ws.Range("M7").Value = -1097.2354
$ws.Range("N7").Value = -2193.5714
$ws.Range("H16").Value = 1927.0625
$ws.Range("I16").Value = 2018.1
$ws.Range("J16").Value = 1775.3334
$ws.Range("K16").Value = 2018.1
$ws.Range("L16").Value = 1775.3334
$ws.Range("M16").Value = -1848.1
$ws.Range("N16").Value = -2115.3334
$ws.Range("H40").Value = 1999.2727
$ws.Range("I40").Value = 1499
$ws.Range("J40").Value = 3333.3333
$ws.Range("K40").Value = 1499
$ws.Range("L40").Value = 3333.3333
$ws.Range("M40").Value = -1363
$ws.Range("N40").Value = -3605.3333
$ws.Range("H126").Value = 1431
$ws.Range("I126").Value = 1209.2354
$ws.Range("J126").Value = 1969.5714
$ws.Range("K126").Value = 3627.7062
$ws.Range("L126").Value = 5908.7142
$ws.Range("M126").Value = -1157.7062
$ws.Range("N126").Value = -10848.7142

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H136").Value = 4008.55
$ws.Range("I136").Value = 1293.8422
$ws.Range("J136").Value = 6464.7144
$ws.Range("K136").Value = 3881.5266
$ws.Range("L136").Value = 19394.1432
$ws.Range("M136").Value = -1331.5266
$ws.Range("N136").Value = -24494.1432
